$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold mixed numeric/text price strings (e.g. "30.317.08",
# "1.000"). Force each target cell to Text format before writing so Excel
# does not coerce numeric-looking values into actual numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.317.08"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.33"
$ws.Range("E3").Value = "  -3.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.87"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7266"
$ws.Range("E6").Value = "  -9.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3287"
$ws.Range("E8").Value = "  -9.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.31"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06819"
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8073"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08052"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.933.28"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.419"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.90"
$ws.Range("E15").Value = "  -5.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.326.33"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.20"
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008005"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.842"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.189.69"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.879"
$ws.Range("E24").Value = "  -4.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.682"
$ws.Range("E25").Value = "  -4.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.21"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.409"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.13"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1334"
$ws.Range("E29").Value = "  -12.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.555"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.337"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.401"
$ws.Range("E32").Value = "  -4.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.188"
$ws.Range("E33").Value = "  -4.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05083"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.221"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7393"
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.751"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.825"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.611"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "79.43"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.991"
$ws.Range("E43").Value = "  -9.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8343"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.32"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.770"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.288"
$ws.Range("E48").Value = "  -4.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.49"
$ws.Range("E49").Value = "  -1.11%  "

# Rows 50 and 51 swap: Cronos <-> NEARProtocol (with updated price/volume)
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.479"
$ws.Range("E50").Value = "  +0.53%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05940"
$ws.Range("E51").Value = "  -0.28%  "
